$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 'WARNING: Sen slope based on tied non-censored values'
$ws.Range("F2").Value = 0.448088767169242
$ws.Range("G2").Value = 0.0188679245283019
$ws.Range("H2").Value = 0.716981132075472
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = -0.687210303096641
$ws.Range("M2").Value = 1.11352781817431
$ws.Range("N2").Value = 0
$ws.Range("P2").Value = 'As likely as not improving'
# Row 3
$ws.Range("H3").Value = 0.8545454545454541
$ws.Range("J3").Value = 10.87
$ws.Range("K3").Value = -0.0322802845739173
$ws.Range("L3").Value = -0.130104226207692
$ws.Range("M3").Value = 0.0544771793234744
$ws.Range("N3").Value = -0.296966739410463
# Row 4
$ws.Range("F4").Value = 0.203401914850544
$ws.Range("G4").Value = 0.196428571428571
$ws.Range("H4").Value = 0.214285714285714
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0.0008877851540782001
# Row 5
$ws.Range("D5").Value = $true
$ws.Range("F5").Value = 0.962008810711416
$ws.Range("H5").Value = 0.660714285714286
$ws.Range("J5").Value = 53.5
$ws.Range("K5").Value = -8.86051212938005
$ws.Range("L5").Value = -13.7145073262018
$ws.Range("M5").Value = -0.838981653601466
$ws.Range("N5").Value = -16.5617049147291
$ws.Range("P5").Value = 'Extremely likely improving'
# Row 6
$ws.Range("E6").Value = '< 5 Non-censored values'
$ws.Range("G6").Value = 0.923076923076923
$ws.Range("H6").Value = 0.0769230769230769
$ws.Range("I6").Value = 1
# Row 7
$ws.Range("F7").Value = 0.991657020460483
$ws.Range("G7").Value = 0.714285714285714
$ws.Range("P7").Value = 'Virtually certain improving'
# Row 8
$ws.Range("E8").Value = 'ok'
$ws.Range("F8").Value = 0.350681023734419
$ws.Range("J8").Value = 0.189
$ws.Range("K8").Value = 0.0072400888585099
$ws.Range("L8").Value = -0.0138846484847568
$ws.Range("M8").Value = 0.0221344276700109
$ws.Range("N8").Value = 3.83073484577244
# Row 9
$ws.Range("F9").Value = 0.233831788510869
$ws.Range("H9").Value = 0.745454545454545
$ws.Range("J9").Value = 7.54
$ws.Range("K9").Value = -0.0208714285714285
$ws.Range("L9").Value = -0.0664168974861666
$ws.Range("M9").Value = 0.0296597021670357
$ws.Range("N9").Value = -0.276809397499052
$ws.Range("P9").Value = 'Unlikely increasing'
# Row 10
$ws.Range("F10").Value = 0.295676810411018
$ws.Range("H10").Value = 0.928571428571429
$ws.Range("J10").Value = 0.19525
$ws.Range("K10").Value = 0.0074897470950102
$ws.Range("L10").Value = -0.0146072817879606
$ws.Range("M10").Value = 0.0209071243374856
$ws.Range("N10").Value = 3.83597802561345
$ws.Range("P10").Value = 'Unlikely improving'
# Row 11
$ws.Range("E11").Value = 'WARNING: Sen slope based on tied non-censored values'
$ws.Range("F11").Value = 0.5
$ws.Range("H11").Value = 0.642857142857143
$ws.Range("J11").Value = 0.255
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = -0.0135578353432446
$ws.Range("M11").Value = 0.0225701422743446
$ws.Range("N11").Value = 0
$ws.Range("P11").Value = 'As likely as not improving'
# Row 12
$ws.Range("E12").Value = 'WARNING: Sen slope based on tied non-censored values'
$ws.Range("F12").Value = 0.265908035745463
$ws.Range("H12").Value = 0.339285714285714
$ws.Range("J12").Value = 0.0105
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = -0.0003618520901061
$ws.Range("M12").Value = 0.0010182682982659
$ws.Range("N12").Value = 0
$ws.Range("P12").Value = 'Unlikely improving'
# Row 13
$ws.Range("F13").Value = 0.0505208312897955
$ws.Range("G13").Value = 0.0092592592592592
$ws.Range("H13").Value = 0.731481481481482
$ws.Range("J13").Value = 5.25
$ws.Range("K13").Value = 0.163569189431258
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0.601873106292202
$ws.Range("N13").Value = 3.11560360821445
$ws.Range("P13").Value = 'Very unlikely improving'
# Row 14
$ws.Range("F14").Value = 0.684092697367081
$ws.Range("H14").Value = 0.839285714285714
$ws.Range("J14").Value = 10.725
$ws.Range("K14").Value = 0.0114635172764673
$ws.Range("L14").Value = -0.0213578637470712
$ws.Range("M14").Value = 0.0445415434048717
$ws.Range("N14").Value = 0.106885941971723
# Row 15
$ws.Range("E15").Value = 'ok'
$ws.Range("F15").Value = 0.997924245435666
$ws.Range("G15").Value = 0.120689655172414
$ws.Range("K15").Value = -0.0003162337662337
$ws.Range("L15").Value = -0.0004993164730006
$ws.Range("N15").Value = -3.51370851370851
$ws.Range("P15").Value = 'Virtually certain improving'
# Row 16
$ws.Range("F16").Value = 0.837190301088992
$ws.Range("G16").Value = 0.916666666666667
$ws.Range("H16").Value = 0.0925925925925926
# Row 17
$ws.Range("F17").Value = 0.999883356486876
$ws.Range("G17").Value = 0.568965517241379
# Row 18
$ws.Range("E18").Value = 'ok'
$ws.Range("F18").Value = 0.78305319351
$ws.Range("H18").Value = 0.9655172413793101
$ws.Range("J18").Value = 0.1914
$ws.Range("K18").Value = -0.0023029774596405
$ws.Range("L18").Value = -0.008559967975134201
$ws.Range("M18").Value = 0.0033229787887986
$ws.Range("N18").Value = -1.20322751287386
$ws.Range("P18").Value = 'Likely improving'
# Row 19
$ws.Range("F19").Value = 0.18063233558513
$ws.Range("H19").Value = 0.616071428571429
$ws.Range("J19").Value = 7.545
$ws.Range("K19").Value = -0.009204117541223
$ws.Range("L19").Value = -0.0248469387755102
$ws.Range("M19").Value = 0.0074814259700623
$ws.Range("N19").Value = -0.121989629439669
$ws.Range("P19").Value = 'Unlikely increasing'
# Row 20
$ws.Range("F20").Value = 0.758613259110031
$ws.Range("H20").Value = 0.913793103448276
$ws.Range("J20").Value = 0.19565
$ws.Range("K20").Value = -0.0023264331210191
$ws.Range("L20").Value = -0.0091842480067456
$ws.Range("M20").Value = 0.003439901482267
$ws.Range("N20").Value = -1.18907902939898
$ws.Range("P20").Value = 'Likely improving'
# Row 21
$ws.Range("D21").Value = $true
$ws.Range("F21").Value = 0.0013594556714927
$ws.Range("G21").Value = 0.0061349693251533
$ws.Range("H21").Value = 0.736196319018405
$ws.Range("J21").Value = 2.8
$ws.Range("K21").Value = 0.115954437361478
$ws.Range("L21").Value = 0.0183066002978353
$ws.Range("M21").Value = 0.322192349437707
$ws.Range("N21").Value = 4.14122990576708
# Row 22
$ws.Range("J22").Value = 0.363
$ws.Range("K22").Value = -0.045292442104493
$ws.Range("M22").Value = 0.0637427676594442
$ws.Range("N22").Value = -12.4772567780972
# Row 23
$ws.Range("J23").Value = 110
$ws.Range("K23").Value = -4.8025053304904
$ws.Range("M23").Value = 1.91611989466209
$ws.Range("N23").Value = -4.36591393680946
# Row 24
$ws.Range("F24").Value = 0.231216363225238
$ws.Range("K24").Value = -0.536112192139103
$ws.Range("M24").Value = 0.496152602534566
$ws.Range("N24").Value = -9.747494402529149
$ws.Range("P24").Value = 'Unlikely improving'
